# Auto-generated edit script: refresh Leve profit-calculation columns
# (currentAveragePrice / LevePrice / LeveProfit, columns H-N) across all
# job sheets per the scheduled runner's updated market-price snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9444.143
$ws.Range("I74").Value = 9444.143
$ws.Range("K74").Value = 9444.143
$ws.Range("M74").Value = -8508.143
$ws.Range("H77").Value = 9444.143
$ws.Range("I77").Value = 9444.143
$ws.Range("K77").Value = 47220.715
$ws.Range("M77").Value = -42540.715
$ws.Range("H132").Value = 4969.5347
$ws.Range("I132").Value = 3137.1428
$ws.Range("J132").Value = 12986.25
$ws.Range("K132").Value = 9411.428400000001
$ws.Range("L132").Value = 38958.75
$ws.Range("M132").Value = -6881.428400000001
$ws.Range("N132").Value = -44018.75
$ws.Range("H135").Value = 6220.6665
$ws.Range("I135").Value = 4907.6665
$ws.Range("J135").Value = 9503.166999999999
$ws.Range("K135").Value = 44168.9985
$ws.Range("L135").Value = 85528.503
$ws.Range("M135").Value = -41633.9985
$ws.Range("N135").Value = -90598.503
$ws.Range("H137").Value = 4219.2324
$ws.Range("I137").Value = 3896.7727
$ws.Range("J137").Value = 4557.048
$ws.Range("K137").Value = 11690.3181
$ws.Range("L137").Value = 13671.144
$ws.Range("M137").Value = -9140.3181
$ws.Range("N137").Value = -18771.144
$ws.Range("H138").Value = 8276.473
$ws.Range("I138").Value = 11947.667
$ws.Range("J138").Value = 7942.727
$ws.Range("K138").Value = 35843.001
$ws.Range("L138").Value = 23828.181
$ws.Range("M138").Value = -30703.001
$ws.Range("N138").Value = -34108.181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3126.3462
$ws.Range("I2").Value = 2905.15
$ws.Range("J2").Value = 3863.6667
$ws.Range("K2").Value = 2905.15
$ws.Range("L2").Value = 3863.6667
$ws.Range("M2").Value = -2792.15
$ws.Range("N2").Value = -4089.6667
$ws.Range("H32").Value = 3237
$ws.Range("I32").Value = 3207.6167
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 3207.6167
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -2920.6167
$ws.Range("N32").Value = -5574
$ws.Range("H109").Value = 20599
$ws.Range("J109").Value = 20599
$ws.Range("L109").Value = 20599
$ws.Range("N109").Value = -23373
$ws.Range("H116").Value = 3126.3462
$ws.Range("I116").Value = 2905.15
$ws.Range("J116").Value = 3863.6667
$ws.Range("K116").Value = 2905.15
$ws.Range("L116").Value = 3863.6667
$ws.Range("M116").Value = -611.1500000000001
$ws.Range("N116").Value = -8451.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3126.3462
$ws.Range("I3").Value = 2905.15
$ws.Range("J3").Value = 3863.6667
$ws.Range("K3").Value = 2905.15
$ws.Range("L3").Value = 3863.6667
$ws.Range("M3").Value = -2791.15
$ws.Range("N3").Value = -4091.6667
$ws.Range("H20").Value = 4923.4165
$ws.Range("I20").Value = 6863.9375
$ws.Range("J20").Value = 1042.375
$ws.Range("K20").Value = 6863.9375
$ws.Range("L20").Value = 1042.375
$ws.Range("M20").Value = -6616.9375
$ws.Range("N20").Value = -1536.375
$ws.Range("H80").Value = 2091.1667
$ws.Range("I80").Value = 1948.1666
$ws.Range("J80").Value = 2234.1667
$ws.Range("K80").Value = 1948.1666
$ws.Range("L80").Value = 2234.1667
$ws.Range("M80").Value = -950.1666
$ws.Range("N80").Value = -4230.1667
$ws.Range("H83").Value = 2091.1667
$ws.Range("I83").Value = 1948.1666
$ws.Range("J83").Value = 2234.1667
$ws.Range("K83").Value = 9740.833000000001
$ws.Range("L83").Value = 11170.8335
$ws.Range("M83").Value = -4748.833000000001
$ws.Range("N83").Value = -21154.8335
$ws.Range("H132").Value = 188333.33
$ws.Range("J132").Value = 188333.33
$ws.Range("L132").Value = 188333.33
$ws.Range("N132").Value = -198453.33
$ws.Range("H134").Value = 10002326
$ws.Range("I134").Value = 2583.889
$ws.Range("J134").Value = 100000000
$ws.Range("K134").Value = 7751.667
$ws.Range("L134").Value = 300000000
$ws.Range("M134").Value = -5216.667
$ws.Range("N134").Value = -300005070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2471.923
$ws.Range("I58").Value = 2248.6
$ws.Range("J58").Value = 3216.3333
$ws.Range("K58").Value = 2248.6
$ws.Range("L58").Value = 3216.3333
$ws.Range("M58").Value = -2045.6
$ws.Range("N58").Value = -3622.3333
$ws.Range("H132").Value = 1883.75
$ws.Range("I132").Value = 1962.5714
$ws.Range("J132").Value = 1699.8334
$ws.Range("K132").Value = 5887.7142
$ws.Range("L132").Value = 5099.5002
$ws.Range("M132").Value = -3357.7142
$ws.Range("N132").Value = -10159.5002
$ws.Range("H136").Value = 2471.923
$ws.Range("I136").Value = 2248.6
$ws.Range("J136").Value = 3216.3333
$ws.Range("K136").Value = 6745.799999999999
$ws.Range("L136").Value = 9648.999899999999
$ws.Range("M136").Value = -4195.799999999999
$ws.Range("N136").Value = -14748.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 4672.1665
$ws.Range("I117").Value = 583
$ws.Range("J117").Value = 5490
$ws.Range("K117").Value = 1749
$ws.Range("L117").Value = 16470
$ws.Range("M117").Value = 1693
$ws.Range("N117").Value = -23354

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8041.222
$ws.Range("I70").Value = 7641.1
$ws.Range("J70").Value = 8276.588
$ws.Range("K70").Value = 7641.1
$ws.Range("L70").Value = 8276.588
$ws.Range("M70").Value = -7371.1
$ws.Range("N70").Value = -8816.588
$ws.Range("H73").Value = 8041.222
$ws.Range("I73").Value = 7641.1
$ws.Range("J73").Value = 8276.588
$ws.Range("K73").Value = 7641.1
$ws.Range("L73").Value = 8276.588
$ws.Range("M73").Value = -6705.1
$ws.Range("N73").Value = -10148.588
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("H132").Value = 6160749
$ws.Range("I132").Value = 2187.2273
$ws.Range("J132").Value = 21215012
$ws.Range("K132").Value = 6561.6819
$ws.Range("L132").Value = 63645036
$ws.Range("M132").Value = -4031.6819
$ws.Range("N132").Value = -63650096
$ws.Range("N104").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14671527
$ws.Range("I22").Value = 14671527
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 14671527
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -14671232
$ws.Range("H23").Value = 4875
$ws.Range("I23").Value = 4875
$ws.Range("K23").Value = 4875
$ws.Range("M23").Value = -4645
$ws.Range("H27").Value = 14671527
$ws.Range("I27").Value = 14671527
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 14671527
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -14671420
$ws.Range("H132").Value = 2912.5715
$ws.Range("I132").Value = 1790.3125
$ws.Range("J132").Value = 4408.9165
$ws.Range("K132").Value = 5370.9375
$ws.Range("L132").Value = 13226.7495
$ws.Range("M132").Value = -2840.9375
$ws.Range("N132").Value = -18286.7495
$ws.Range("H136").Value = 2469.8823
$ws.Range("I136").Value = 2418
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 7254
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -4704
$ws.Range("N136").Value = -15000
$ws.Range("N22").ClearContents()
$ws.Range("N27").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1021.55554
$ws.Range("I113").Value = 741.5714
$ws.Range("J113").Value = 2001.5
$ws.Range("K113").Value = 2224.7142
$ws.Range("L113").Value = 6004.5
$ws.Range("M113").Value = -54.71420000000035
$ws.Range("N113").Value = -10344.5
$ws.Range("H122").Value = 4637.25
$ws.Range("I122").Value = 4275
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 12825
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -10375
$ws.Range("N122").Value = -19898.5
$ws.Range("H136").Value = 319261.7
$ws.Range("I136").Value = 8261.5
$ws.Range("J136").Value = 1666929.1
$ws.Range("K136").Value = 24784.5
$ws.Range("L136").Value = 5000787.300000001
$ws.Range("M136").Value = -22234.5
$ws.Range("N136").Value = -5005887.300000001
